# Generate Report for Handoff
# Adds a new tracked file (d2439d33-...) as "Ready for handoff" to the
# Overview, zh-cn and de-de localization-status sheets, mirroring the
# existing 752d6c67-... row.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# ---- shared literal values -------------------------------------------------

$fileName   = "d2439d33-6f28-4230-9afb-8d36418dcea2ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$pathName   = "e2e\d2439d33-6f28-4230-9afb-8d36418dcea2ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$extension  = ".md"
$status     = "Ready for handoff"
$genDate    = "2016-09-06 18:37:42"

$zhHandoffFile = "d2439d33-6f28-4230-9afb-8d36418dcea2oooooooooooooooooooooooooooooooooooooooo.ed59ca03e236a28a956550b8ca7a552572793ce4.zh-cn.xlf"
$zhHandoffDate = "2016-09-06 18:37:37"

$deHandoffFile = "d2439d33-6f28-4230-9afb-8d36418dcea2oooooooooooooooooooooooooooooooooooooooo.ed59ca03e236a28a956550b8ca7a552572793ce4.de-de.xlf"
$deHandoffDate = "2016-09-06 18:37:42"

$sourcePath = "e2e"
$priority   = "ht"
$dupFlag    = "'False"
$toLocalize = "'True"
$hasMeta    = "'False"
$backDate   = "0001-01-01 00:00:00"
$emptyCell  = "'"

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1e65ebc951d995060ca80a9696d109c730db3941/e2e/" + $fileName

# ---- Overview sheet ---------------------------------------------------------

$lo1 = $ws1.ListObjects.Item(1)
$lo1.ListRows.Add() | Out-Null

$ws1.Cells.Item(3,1).Value = $fileName                  # A3 File Name
$ws1.Cells.Item(3,2).Value = $pathName                   # B3 Path And Name
$ws1.Cells.Item(3,3).Value = $extension                  # C3 Extension
$ws1.Cells.Item(3,4).Value = $emptyCell                   # D3 Publish URL
$ws1.Cells.Item(3,5).Value = $status                      # E3 zh-cn
$ws1.Cells.Item(3,6).Value = $status                      # F3 de-de
$ws1.Cells.Item(3,7).Value = $genDate                     # G3 Latest HO Xliff Generate Date
$ws1.Cells.Item(3,7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws1.Hyperlinks.Add($ws1.Range("B3"), $baseUrl, "", "", $pathName) | Out-Null

# ---- zh-cn sheet -------------------------------------------------------------

$lo2 = $ws2.ListObjects.Item(1)
$lo2.ListRows.Add() | Out-Null

$ws2.Cells.Item(3,1).Value  = $fileName        # A3 Source File Name
$ws2.Cells.Item(3,2).Value  = $extension       # B3 File Extension
$ws2.Cells.Item(3,3).Value  = $status          # C3 Status
$ws2.Cells.Item(3,4).Value  = $sourcePath      # D3 Source Path
$ws2.Cells.Item(3,5).Value  = $priority        # E3 Priority
$ws2.Cells.Item(3,6).Value  = $dupFlag         # F3 Content Duplicate
$ws2.Cells.Item(3,7).Value  = $zhHandoffFile   # G3 Latest Handoff File
$ws2.Cells.Item(3,8).Value  = $zhHandoffDate   # H3 Latest Handoff Datetime
$ws2.Cells.Item(3,9).Value  = $emptyCell        # I3 Latest Target File
$ws2.Cells.Item(3,10).Value = $emptyCell        # J3 Latest Handback File
$ws2.Cells.Item(3,11).Value = $backDate        # K3 Latest Handback DateTime
$ws2.Cells.Item(3,12).Value = $emptyCell        # L3 Reference Tokens
$ws2.Cells.Item(3,13).Value = $toLocalize      # M3 To be localized
$ws2.Cells.Item(3,14).Value = $emptyCell        # N3 Dependency From
$ws2.Cells.Item(3,15).Value = $hasMeta         # O3 Has metadata
$ws2.Cells.Item(3,16).Value = $emptyCell        # P3 Error Detail

$ws2.Cells.Item(3,8).NumberFormat  = "yyyy-mm-dd HH:mm:ss"
$ws2.Cells.Item(3,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws2.Hyperlinks.Add($ws2.Range("A3"), $baseUrl, "", "", $fileName) | Out-Null

# ---- de-de sheet -------------------------------------------------------------

$lo3 = $ws3.ListObjects.Item(1)
$lo3.ListRows.Add() | Out-Null

$ws3.Cells.Item(3,1).Value  = $fileName        # A3 Source File Name
$ws3.Cells.Item(3,2).Value  = $extension       # B3 File Extension
$ws3.Cells.Item(3,3).Value  = $status          # C3 Status
$ws3.Cells.Item(3,4).Value  = $sourcePath      # D3 Source Path
$ws3.Cells.Item(3,5).Value  = $priority        # E3 Priority
$ws3.Cells.Item(3,6).Value  = $dupFlag         # F3 Content Duplicate
$ws3.Cells.Item(3,7).Value  = $deHandoffFile   # G3 Latest Handoff File
$ws3.Cells.Item(3,8).Value  = $genDate         # H3 Latest Handoff Datetime
$ws3.Cells.Item(3,9).Value  = $emptyCell        # I3 Latest Target File
$ws3.Cells.Item(3,10).Value = $emptyCell        # J3 Latest Handback File
$ws3.Cells.Item(3,11).Value = $backDate        # K3 Latest Handback DateTime
$ws3.Cells.Item(3,12).Value = $emptyCell        # L3 Reference Tokens
$ws3.Cells.Item(3,13).Value = $toLocalize      # M3 To be localized
$ws3.Cells.Item(3,14).Value = $emptyCell        # N3 Dependency From
$ws3.Cells.Item(3,15).Value = $hasMeta         # O3 Has metadata
$ws3.Cells.Item(3,16).Value = $emptyCell        # P3 Error Detail

$ws3.Cells.Item(3,8).NumberFormat  = "yyyy-mm-dd HH:mm:ss"
$ws3.Cells.Item(3,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws3.Hyperlinks.Add($ws3.Range("A3"), $baseUrl, "", "", $fileName) | Out-Null
